$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.150936007499695
$ws.Range("B1").Value = 2.135812044143677
$ws.Range("C1").Value = 10.23412227630615
$ws.Range("D1").Value = 2.542796850204468
$ws.Range("E1").Value = 1.274287939071655
